# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The two "Periodo Mora" / "Valor Mora" rows (16 and 17) need to swap places:
#   Row 16 currently shows Periodo 2503 / Valor 56940  -> should become Periodo 2502 / Valor 47450
#   Row 17 currently shows Periodo 2502 / Valor 47450   -> should become Periodo 2503 / Valor 56940
# i.e. the (Periodo, Valor) pairs stay matched together, just the two rows trade places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current values before overwriting anything.
$periodo16 = $ws.Range("E16").Value2
$periodo17 = $ws.Range("E17").Value2
$valor16   = $ws.Range("F16").Value2
$valor17   = $ws.Range("F17").Value2

# Swap row 16 <-> row 17 for both the period label and the overdue value.
$ws.Range("E16").Value = $periodo17
$ws.Range("E17").Value = $periodo16
$ws.Range("F16").Value = $valor17
$ws.Range("F17").Value = $valor16
